# Update "想去人数" (want-to-go count) figures refreshed at commit 456a3b4.
# Two sheets carry the same event rows and need the same F-column updates:
#   展览 (Exhibitions)   - rows 2-21
#   全部类型 (All types)  - rows 2-23 (superset incl. 演出 sheet rows)

$wb = $excel.ActiveWorkbook

$exhibitions = $wb.Worksheets.Item("展览")
$exhibitions.Range("F2").Value = 1203
$exhibitions.Range("F3").Value = 428
$exhibitions.Range("F4").Value = 287
$exhibitions.Range("F5").Value = 150
$exhibitions.Range("F7").Value = 12396
$exhibitions.Range("F10").Value = 19
$exhibitions.Range("F11").Value = 1
$exhibitions.Range("F13").Value = 12227
$exhibitions.Range("F14").Value = 4852
$exhibitions.Range("F15").Value = 4728
$exhibitions.Range("F16").Value = 139
$exhibitions.Range("F20").Value = 953
$exhibitions.Range("F21").Value = 5

$allTypes = $wb.Worksheets.Item("全部类型")
$allTypes.Range("F2").Value = 1203
$allTypes.Range("F3").Value = 428
$allTypes.Range("F4").Value = 287
$allTypes.Range("F5").Value = 150
$allTypes.Range("F9").Value = 12396
$allTypes.Range("F12").Value = 19
$allTypes.Range("F13").Value = 1
$allTypes.Range("F15").Value = 12227
$allTypes.Range("F16").Value = 4852
$allTypes.Range("F17").Value = 4728
$allTypes.Range("F18").Value = 139
$allTypes.Range("F22").Value = 953
$allTypes.Range("F23").Value = 5
